# Actualización 10 de Mayo
# Updates the D/E/F/G/H statistics columns on the three "Estadisticos"
# sheets (1P, 2P, Final). "Blancos" (D) goes to 0, "Aprobados" (F) is
# raised to match "Totales" (C), "Por_Apro" (G) becomes 100, and the
# "Promedio" (H) column is refreshed with new values (added outright on
# the 2P sheet, which didn't have an H column populated yet).

$wb = $excel.ActiveWorkbook

# Row data: Grupo row number -> D, F, G, H values common to all three sheets.
$rows = @(
    @{ Row = 2; D = 0; F = 36; G = 100; H = 8.1 },
    @{ Row = 3; D = 0; F = 36; G = 100; H = 8.6 },
    @{ Row = 4; D = 0; F = 35; G = 100; H = 7.9 },
    @{ Row = 5; D = 0; F = 25; G = 100; H = 8.4 }
)

foreach ($sheetName in @("Estadisticos 1P", "Estadisticos 2P")) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($r in $rows) {
        $ws.Cells.Item($r.Row, 4).Value = $r.D   # D: Blancos
        $ws.Cells.Item($r.Row, 5).Value = 0       # E: Reprobados
        $ws.Cells.Item($r.Row, 6).Value = $r.F   # F: Aprobados
        $ws.Cells.Item($r.Row, 7).Value = $r.G   # G: Por_Apro
        $ws.Cells.Item($r.Row, 8).Value = $r.H   # H: Promedio
    }
}

# "Estadisticos Final" matches the others except rows 3 and 5 keep their
# original Promedio (H) values.
$wsFinal = $wb.Worksheets.Item("Estadisticos Final")
$wsFinal.Cells.Item(2, 4).Value = 0
$wsFinal.Cells.Item(2, 5).Value = 0
$wsFinal.Cells.Item(2, 6).Value = 36
$wsFinal.Cells.Item(2, 7).Value = 100
$wsFinal.Cells.Item(2, 8).Value = 8.1

$wsFinal.Cells.Item(3, 4).Value = 0
$wsFinal.Cells.Item(3, 5).Value = 0
$wsFinal.Cells.Item(3, 6).Value = 36
$wsFinal.Cells.Item(3, 7).Value = 100

$wsFinal.Cells.Item(4, 4).Value = 0
$wsFinal.Cells.Item(4, 5).Value = 0
$wsFinal.Cells.Item(4, 6).Value = 35
$wsFinal.Cells.Item(4, 7).Value = 100
$wsFinal.Cells.Item(4, 8).Value = 7.9

$wsFinal.Cells.Item(5, 4).Value = 0
$wsFinal.Cells.Item(5, 5).Value = 0
$wsFinal.Cells.Item(5, 6).Value = 25
$wsFinal.Cells.Item(5, 7).Value = 100
